$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.56338
$ws.Range("H2").Value = 1.69014
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.774269
$ws.Range("N2").Value = 2.322807
$ws.Range("O2").Value = 0.02746805195987118
$ws.Range("P2").Value = 0.02746805195987118
$ws.Range("Q2").Value = 0.43620766922
$ws.Range("R2").Value = 3.92586902298
$ws.Range("S2").Value = 0.02746805195987118
$ws.Range("T2").Value = 0.02746805195987118

# Row 3
$ws.Range("G3").Value = 0.56338
$ws.Range("H3").Value = 1.69014
$ws.Range("M3").Value = 25.63013966666666
$ws.Range("N3").Value = 76.89041899999999
$ws.Range("O3").Value = 0.9092576457313354
$ws.Range("P3").Value = 0.9092576457313354
$ws.Range("Q3").Value = 14.43950808540666
$ws.Range("R3").Value = 129.95557276866
$ws.Range("S3").Value = 0.9092576457313354
$ws.Range("T3").Value = 0.9092576457313354

# Row 4
$ws.Range("G4").Value = 0.56338
$ws.Range("H4").Value = 1.69014
$ws.Range("M4").Value = 1.783575
$ws.Range("N4").Value = 5.350725000000001
$ws.Range("O4").Value = 0.06327430230879351
$ws.Range("P4").Value = 0.06327430230879351
$ws.Range("Q4").Value = 1.0048304835
$ws.Range("R4").Value = 9.0434743515
$ws.Range("S4").Value = 0.06327430230879351
$ws.Range("T4").Value = 0.06327430230879351
